$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 420036.03
$ws.Range("J17").Value = 436559.34
$ws.Range("L17").Value = 1309678.02
$ws.Range("N17").Value = -1310014.02
$ws.Range("H51").Value = 9800.200000000001
$ws.Range("I51").Value = 9600.4
$ws.Range("J51").Value = 10000
$ws.Range("K51").Value = 9600.4
$ws.Range("L51").Value = 10000
$ws.Range("M51").Value = -9116.4
$ws.Range("N51").Value = -10968
$ws.Range("H70").Value = 948510.4
$ws.Range("I70").Value = 1384084.4
$ws.Range("K70").Value = 4152253.2
$ws.Range("M70").Value = -4151983.2
$ws.Range("H73").Value = 948510.4
$ws.Range("I73").Value = 1384084.4
$ws.Range("K73").Value = 4152253.2
$ws.Range("M73").Value = -4151317.2
$ws.Range("H101").Value = 416
$ws.Range("I101").Value = 416
$ws.Range("K101").Value = 1248
$ws.Range("M101").Value = 374
$ws.Range("H138").Value = 4202.5884
$ws.Range("I138").Value = 3855.7273
$ws.Range("J138").Value = 4368.478
$ws.Range("K138").Value = 11567.1819
$ws.Range("L138").Value = 13105.434
$ws.Range("M138").Value = -6427.1819
$ws.Range("N138").Value = -23385.434

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1583.2667
$ws.Range("I2").Value = 1340.909
$ws.Range("K2").Value = 1340.909
$ws.Range("M2").Value = -1227.909
$ws.Range("H74").Value = 4311.577
$ws.Range("I74").Value = 3416.3157
$ws.Range("J74").Value = 6741.5713
$ws.Range("K74").Value = 3416.3157
$ws.Range("L74").Value = 6741.5713
$ws.Range("M74").Value = -2542.3157
$ws.Range("N74").Value = -8489.5713
$ws.Range("H77").Value = 4311.577
$ws.Range("I77").Value = 3416.3157
$ws.Range("J77").Value = 6741.5713
$ws.Range("K77").Value = 17081.5785
$ws.Range("L77").Value = 33707.85649999999
$ws.Range("M77").Value = -12713.5785
$ws.Range("N77").Value = -42443.85649999999
$ws.Range("H116").Value = 1583.2667
$ws.Range("I116").Value = 1340.909
$ws.Range("K116").Value = 1340.909
$ws.Range("M116").Value = 953.0909999999999
$ws.Range("H132").Value = 4123.731
$ws.Range("I132").Value = 3671.9092
$ws.Range("J132").Value = 6608.75
$ws.Range("K132").Value = 11015.7276
$ws.Range("L132").Value = 19826.25
$ws.Range("M132").Value = -8485.7276
$ws.Range("N132").Value = -24886.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1583.2667
$ws.Range("I3").Value = 1340.909
$ws.Range("K3").Value = 1340.909
$ws.Range("M3").Value = -1226.909
$ws.Range("H82").Value = 22716.666
$ws.Range("I82").Value = 13660
$ws.Range("J82").Value = 68000
$ws.Range("K82").Value = 13660
$ws.Range("L82").Value = 68000
$ws.Range("M82").Value = -13277
$ws.Range("N82").Value = -68766
$ws.Range("H85").Value = 22716.666
$ws.Range("I85").Value = 13660
$ws.Range("J85").Value = 68000
$ws.Range("K85").Value = 13660
$ws.Range("L85").Value = 68000
$ws.Range("M85").Value = -12334
$ws.Range("N85").Value = -70652
$ws.Range("H86").Value = 5119.3
$ws.Range("I86").Value = 2866.5
$ws.Range("K86").Value = 2866.5
$ws.Range("M86").Value = -1743.5
$ws.Range("H89").Value = 5119.3
$ws.Range("I89").Value = 2866.5
$ws.Range("K89").Value = 14332.5
$ws.Range("M89").Value = -8716.5
$ws.Range("H94").Value = 1509.3793
$ws.Range("I94").Value = 1472.4
$ws.Range("K94").Value = 1472.4
$ws.Range("M94").Value = -1021.4
$ws.Range("H105").Value = 2227.2
$ws.Range("I105").Value = 2100.158
$ws.Range("K105").Value = 2100.158
$ws.Range("M105").Value = -353.1579999999999
$ws.Range("H107").Value = 41836016
$ws.Range("I107").Value = 288434.72
$ws.Range("K107").Value = 288434.72
$ws.Range("M107").Value = -286514.72

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4300.029
$ws.Range("I31").Value = 2774.25
$ws.Range("K31").Value = 2774.25
$ws.Range("M31").Value = -2479.25
$ws.Range("H34").Value = 4300.029
$ws.Range("I34").Value = 2774.25
$ws.Range("K34").Value = 2774.25
$ws.Range("M34").Value = -2572.25
$ws.Range("H62").Value = 170645
$ws.Range("I62").Value = 5126.25
$ws.Range("J62").Value = 303060
$ws.Range("K62").Value = 5126.25
$ws.Range("L62").Value = 303060
$ws.Range("M62").Value = -4502.25
$ws.Range("N62").Value = -304308
$ws.Range("H65").Value = 170645
$ws.Range("I65").Value = 5126.25
$ws.Range("J65").Value = 303060
$ws.Range("K65").Value = 25631.25
$ws.Range("L65").Value = 1515300
$ws.Range("M65").Value = -22511.25
$ws.Range("N65").Value = -1521540
$ws.Range("H107").Value = 1135.875
$ws.Range("I107").Value = 1144.3846
$ws.Range("K107").Value = 1144.3846
$ws.Range("M107").Value = 775.6153999999999
$ws.Range("H132").Value = 437654.3
$ws.Range("I132").Value = 2740.1428
$ws.Range("K132").Value = 8220.428400000001
$ws.Range("M132").Value = -5690.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 11877.75
$ws.Range("J81").Value = 13431.714
$ws.Range("L81").Value = 40295.142
$ws.Range("N81").Value = -42541.142
$ws.Range("H84").Value = 11877.75
$ws.Range("J84").Value = 13431.714
$ws.Range("L84").Value = 120885.426
$ws.Range("N84").Value = -132117.426
$ws.Range("H132").Value = 2725.4324
$ws.Range("J132").Value = 2828.4924
$ws.Range("L132").Value = 25456.4316
$ws.Range("N132").Value = -30516.4316
$ws.Range("H137").Value = 3191.4546
$ws.Range("J137").Value = 3643.125
$ws.Range("L137").Value = 10929.375
$ws.Range("N137").Value = -21129.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4275.8184
$ws.Range("I113").Value = 3446.25
$ws.Range("J113").Value = 4749.857
$ws.Range("K113").Value = 3446.25
$ws.Range("L113").Value = 4749.857
$ws.Range("M113").Value = -1276.25
$ws.Range("N113").Value = -9089.857
$ws.Range("H132").Value = 6524.1113
$ws.Range("I132").Value = 5718.684
$ws.Range("K132").Value = 17156.052
$ws.Range("M132").Value = -14626.052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2957.432
$ws.Range("I46").Value = 1748.75
$ws.Range("J46").Value = 3410.6875
$ws.Range("K46").Value = 1748.75
$ws.Range("L46").Value = 3410.6875
$ws.Range("M46").Value = -1560.75
$ws.Range("N46").Value = -3786.6875
$ws.Range("H61").Value = 1234.1666
$ws.Range("I61").Value = 1081
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1081
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -879
$ws.Range("N61").Value = -2404
$ws.Range("H93").Value = 1368.591
$ws.Range("J93").Value = 1750
$ws.Range("L93").Value = 1750
$ws.Range("N93").Value = -4246
$ws.Range("H113").Value = 1234.1666
$ws.Range("I113").Value = 1081
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1081
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1089
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 5459.469
$ws.Range("I132").Value = 3590.4119
$ws.Range("K132").Value = 10771.2357
$ws.Range("M132").Value = -8241.235700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2958.6365
$ws.Range("I96").Value = 1895
$ws.Range("J96").Value = 3065
$ws.Range("K96").Value = 1895
$ws.Range("L96").Value = 3065
$ws.Range("M96").Value = -522
$ws.Range("N96").Value = -5811
$ws.Range("H100").Value = 3447.087
$ws.Range("I100").Value = 3967.6316
$ws.Range("J100").Value = 974.5
$ws.Range("K100").Value = 7935.2632
$ws.Range("L100").Value = 1949
$ws.Range("M100").Value = -7394.2632
$ws.Range("N100").Value = -3031
$ws.Range("H113").Value = 5913.6665
$ws.Range("I113").Value = 6668.778
$ws.Range("J113").Value = 4403.4443
$ws.Range("K113").Value = 20006.334
$ws.Range("L113").Value = 13210.3329
$ws.Range("M113").Value = -17836.334
$ws.Range("N113").Value = -17550.3329
$ws.Range("H132").Value = 466239.2
$ws.Range("I132").Value = 629993.8
$ws.Range("J132").Value = 29560.166
$ws.Range("K132").Value = 1889981.4
$ws.Range("L132").Value = 88680.49800000001
$ws.Range("M132").Value = -1887451.4
$ws.Range("N132").Value = -93740.49800000001
